{"js": "// \"Version 2.\" -> \"Version 1.\"\n// The original text is split across runs (\"Versi\" + \"on\" + \" 2\" + \".\"). We\n// locate the relevant substrings with search() and rewrite them in place so\n// the resulting run layout mirrors the target: \"Version\" then \" 1.\" (the\n// trailing \".\" run disappears, merged into the \" 1.\" run).\n\nconst body = context.document.body;\n\n// \"Versi\" + \"on\" -> merge into a single run \"Version\".\nconst versionRanges = body.search(\"Version\", { matchCase: true });\nawait context.sync();\nversionRanges.items[0].insertText(\"Version\", \"Replace\");\nawait context.sync();\n\n// \" 2\" + \".\" -> merge into a single run \" 1.\".\nconst tailRanges = body.search(\" 2.\", { matchCase: true });\nawait context.sync();\ntailRanges.items[0].insertText(\" 1.\", \"Replace\");\nawait context.sync();\n", "ps1": "# \"Version 2.\" -> \"Version 1.\"\n#\n# The paragraph's text is split across several runs:\n#   \"Versi\" + \"on\" + \" 2\" + (bookmark \"_GoBack\") + \".\"\n# The target layout merges \"Versi\"+\"on\" into one run (\"Version\") and merges\n# \" 2\" with the following \".\" into one run (\" 1.\"), dropping the now-empty\n# trailing run, while keeping the _GoBack bookmark in place between them.\n\n$d = $word.ActiveDocument\n\n# Step 1: merge \"Versi\" + \"on\" into a single run \"Version\".\n$rngVersion = $d.Content\n$null = $rngVersion.Find.Execute(\"Version\", $false, $false, $false, $false, $false, $true, 1, $false, \"Version\", 2)\n\n# Step 2: remove the trailing \".\" run, which sits right after the _GoBack\n# bookmark. Scope the search to the range after the bookmark so only that\n# run's period is touched.\n$bm = $d.Bookmarks.Item(\"_GoBack\")\n$tail = $d.Range($bm.End, $d.Content.End - 1)\n$null = $tail.Find.Execute(\".\", $false, $false, $false, $false, $false, $true, 1, $false, \"\", 2)\n\n# Step 3: turn \" 2\" (the run right before the bookmark) into \" 1.\". Scope the\n# search to the range before the bookmark for the same reason as above.\n$bm = $d.Bookmarks.Item(\"_GoBack\")\n$head = $d.Range(0, $bm.Start)\n$null = $head.Find.Execute(\" 2\", $false, $false, $false, $false, $false, $true, 1, $false, \" 1.\", 2)\n"}
